$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ENGLAND - Premier League: Wolves vs Brighton, SCHEDULED
$ws.Range("D2").Value = "13:00"
$ws.Range("E2").Value = "Wolves"
$ws.Range("F2").Value = "https://www.flashscore.com/res/image/data/OMUzjDkC-rawILjE1.png"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "https://www.flashscore.com/res/image/data/G0q9xjRq-IccPlIAs.png"
$ws.Range("I2").Value = "Brighton"
$ws.Range("J2").Value = "SCHEDULED"

# Row 3 - ENGLAND - Premier League: Aston Villa vs Manchester Utd, SCHEDULED
$ws.Range("D3").Value = "15:05"
$ws.Range("E3").Value = "Aston Villa"
$ws.Range("F3").Value = "https://www.flashscore.com/res/image/data/UHchCEVH-jm1Xyzp7.png"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "https://www.flashscore.com/res/image/data/GhGV3qjT-rBodzytr.png"
$ws.Range("I3").Value = "Manchester Utd"
$ws.Range("J3").Value = "SCHEDULED"

# Row 4 - ENGLAND - Premier League: West Ham vs Everton, SCHEDULED
$ws.Range("D4").Value = "17:30"
$ws.Range("E4").Value = "West Ham"
$ws.Range("F4").Value = "https://www.flashscore.com/res/image/data/YeSfKGlC-hrtlQ906.png"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "https://www.flashscore.com/res/image/data/EBfZuwme-Onr593up.png"
$ws.Range("I4").Value = "Everton"
$ws.Range("J4").Value = "SCHEDULED"

# Row 5 - ENGLAND - Premier League: Arsenal vs West Brom, SCHEDULED (already)
$ws.Range("D5").Value = "20:00"
$ws.Range("E5").Value = "Arsenal"
$ws.Range("F5").Value = "https://www.flashscore.com/res/image/data/0n1ffK6k-pU2IsJm8.png"
$ws.Range("H5").Value = "https://www.flashscore.com/res/image/data/QsGXnZjC-hUScfdXD.png"
$ws.Range("I5").Value = "West Brom"

# Row 6 - FRANCE - Ligue 1: St Etienne vs Marseille, SCHEDULED
$ws.Range("D6").Value = "13:00"
$ws.Range("E6").Value = "St Etienne"
$ws.Range("F6").Value = "https://www.flashscore.com/res/image/data/MF4bIRPq-Qk1sBuEM.png"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = "https://www.flashscore.com/res/image/data/t6Kb5X76-Qk1sBuEM.png"
$ws.Range("I6").Value = "Marseille"
$ws.Range("J6").Value = "SCHEDULED"

# Row 7 - FRANCE - Ligue 1: Angers vs Dijon, SCHEDULED
$ws.Range("D7").Value = "15:00"
$ws.Range("E7").Value = "Angers"
$ws.Range("F7").Value = "https://www.flashscore.com/res/image/data/IN9Ib7jT-EorrQF3M.png"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = "https://www.flashscore.com/res/image/data/IVjOpykC-YHD0XcTg.png"
$ws.Range("I7").Value = "Dijon"
$ws.Range("J7").Value = "SCHEDULED"

# Row 8 - FRANCE - Ligue 1 (was GERMANY - Bundesliga): Metz vs Nimes, SCHEDULED
$ws.Range("C8").Value = "FRANCE - Ligue 1"
$ws.Range("D8").Value = "15:00"
$ws.Range("E8").Value = "Metz"
$ws.Range("F8").Value = "https://www.flashscore.com/res/image/data/SCdF7fjT-4pdbha5J.png"
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = "https://www.flashscore.com/res/image/data/xtTjGRT0-Sr136dyl.png"
$ws.Range("I8").Value = "Nimes"
$ws.Range("J8").Value = "SCHEDULED"

# Row 9 - FRANCE - Ligue 1 (was GERMANY - Bundesliga): Nice vs Brest, SCHEDULED
$ws.Range("C9").Value = "FRANCE - Ligue 1"
$ws.Range("D9").Value = "15:00"
$ws.Range("E9").Value = "Nice"
$ws.Range("F9").Value = "https://www.flashscore.com/res/image/data/bJCGymle-jiROSmFn.png"
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = "https://www.flashscore.com/res/image/data/lWu2w4ne-fXJWG6Mg.png"
$ws.Range("I9").Value = "Brest"
$ws.Range("J9").Value = "SCHEDULED"

# Row 10 - FRANCE - Ligue 1 (was GERMANY - Bundesliga): Strasbourg vs Montpellier, SCHEDULED
$ws.Range("C10").Value = "FRANCE - Ligue 1"
$ws.Range("D10").Value = "15:00"
$ws.Range("E10").Value = "Strasbourg"
$ws.Range("F10").Value = "https://www.flashscore.com/res/image/data/2XOhzSQq-4SOLY7oe.png"
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = "https://www.flashscore.com/res/image/data/U3fC5I96-EkChIdLk.png"
$ws.Range("I10").Value = "Montpellier"
$ws.Range("J10").Value = "SCHEDULED"
